$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.973.56'
$ws.Range("E2").Value = '  -1.69%  '
$ws.Range("D3").Value = '2.169.96'
$ws.Range("E3").Value = '  -2.70%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.88'
$ws.Range("E5").Value = '  -2.02%  '
$ws.Range("E6").Value = '  -2.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '66.11'
$ws.Range("E7").Value = '  -6.86%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.563'
$ws.Range("E9").Value = '  -1.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.88'
$ws.Range("E10").Value = '  -0.34%  '
$ws.Range("E11").Value = '  -4.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '35.51'
$ws.Range("E12").Value = '  -16.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.103'
$ws.Range("E13").Value = '  -2.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.79'
$ws.Range("E14").Value = '  -3.17%  '
$ws.Range("D15").Value = '2.493.37'
$ws.Range("E15").Value = '  -2.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.853'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("E17").Value = '  -4.63%  '
$ws.Range("D18").Value = '2.171.94'
$ws.Range("E18").Value = '  -2.72%  '
$ws.Range("D19").Value = '40.858.31'
$ws.Range("D20").Value = '0.0₃0937'
$ws.Range("E20").Value = '  -2.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.08'
$ws.Range("E21").Value = '  -2.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.21'
$ws.Range("E22").Value = '  -2.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '228.73'
$ws.Range("E23").Value = '  -2.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.09'
$ws.Range("E24").Value = '  -7.85%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.26'
$ws.Range("E26").Value = '  +8.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.68'
$ws.Range("E27").Value = '  -5.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.41'
$ws.Range("E28").Value = '  -3.64%  '
$ws.Range("E29").Value = '  -5.78%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.01'
$ws.Range("E30").Value = '  -1.48%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.00'
$ws.Range("E31").Value = '  -8.91%  '
$ws.Range("E32").Value = '  -2.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.120'
$ws.Range("E33").Value = '  -1.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.66'
$ws.Range("E34").Value = '  +1.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0747'
$ws.Range("E35").Value = '  +3.52%  '
$ws.Range("E36").Value = '  -3.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.55'
$ws.Range("E37").Value = '  -2.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.94'
$ws.Range("E38").Value = '  -1.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.65'
$ws.Range("E39").Value = '  -8.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0305'
$ws.Range("E40").Value = '  +5.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.18'
$ws.Range("E41").Value = '  -5.40%  '
$ws.Range("E42").Value = '  -9.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.36'
$ws.Range("E43").Value = '  -5.55%  '
$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.83'
$ws.Range("E44").Value = '  -5.13%  '
$ws.Range("B45").Value = 'MultiversX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '60.16'
$ws.Range("E45").Value = '  -13.24%  '
$ws.Range("E46").Value = '  -7.76%  '
$ws.Range("E47").Value = '  -0.15%  '
$ws.Range("E48").Value = '  -4.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0987'
$ws.Range("E49").Value = '  -3.19%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.13'
$ws.Range("E50").Value = '  -2.11%  '
$ws.Range("E51").Value = '  -4.21%  '
